$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 216, pushing the existing rows 216:247 down to 217:248.
$ws.Rows(216).Insert()

# Populate the newly inserted row with its data (mirrors the surrounding
# rows' values for the columns that stay constant in this table).
$ws.Cells.Item(216, 1).Value = 1
$ws.Cells.Item(216, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(216, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(216, 4).Value = 45127
$ws.Cells.Item(216, 5).Value = 15
$ws.Cells.Item(216, 6).Value = 100114001
$ws.Cells.Item(216, 7).Value = "Papa"
$ws.Cells.Item(216, 8).Value = "Yagana"
$ws.Cells.Item(216, 9).Value = "1a (cosecha)"
$ws.Cells.Item(216, 10).Value = 1000
$ws.Cells.Item(216, 11).Value = 20000
$ws.Cells.Item(216, 12).Value = 21000
$ws.Cells.Item(216, 13).Value = 20500
$ws.Cells.Item(216, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(216, 15).Value = "Región de Los Lagos"
$ws.Cells.Item(216, 16).Value = 820
$ws.Cells.Item(216, 17).Value = 25
$ws.Cells.Item(216, 18).Value = "Hortaliza"
